$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Save" column - reuse the same header formatting
# (bold font + border + centered alignment) as the existing headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save = 1 when the row's sum (column G) exceeds 8, else 0 — mirrors how
# column F ("Win") already flags notable games. Write literal values (not
# formulas) to match the source data's convention of storing computed
# numbers directly.
for ($r = 2; $r -le 48; $r++) {
    $g = $ws.Range("G$r").Value2
    if ($g -gt 8) {
        $ws.Range("H$r").Value = 1
    } else {
        $ws.Range("H$r").Value = 0
    }
}

$wb.Save()
